$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values for several rows, per repull/mean calculation fix
$ws.Range("F2").Value = 0
$ws.Range("F5").Value = -1
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -2
$ws.Range("F11").Value = -12
$ws.Range("F15").Value = -11
